$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.259.02"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.908.59"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.85"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5255"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3825"
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07315"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.68"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08221"
$ws.Range("E12").Value = "  -3.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.39"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.376"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.677.85"
$ws.Range("E15").Value = "  -11.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008694"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.78"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.291.98"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.83"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.502"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.343"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.16"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.26"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.739"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.99"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.860"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.882"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09259"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8259"
$ws.Range("E32").Value = "  +2.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05081"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.230"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.365"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  +3.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5766"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.085"
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.074"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.616"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.16"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1525"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4935"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.18"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.643"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.66"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06054"
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.03"
$ws.Range("E51").Value = "  -0.05%  "
